$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (matching the workbook's original inline-string cell type for the Price column).
$textCells = @(
    'D5'
    'D6'
    'D7'
    'D9'
    'D10'
    'D11'
    'D12'
    'D14'
    'D19'
    'D21'
    'D22'
    'D23'
    'D24'
    'D25'
    'D26'
    'D27'
    'D28'
    'D29'
    'D30'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D36'
    'D37'
    'D38'
    'D39'
    'D40'
    'D41'
    'D42'
    'D43'
    'D44'
    'D45'
    'D47'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    'D2' = '43.170.90'
    'E2' = '  +1.62%  '
    'D3' = '2.377.97'
    'E3' = '  +6.73%  '
    'D5' = '322.16'
    'E5' = '  +9.06%  '
    'D6' = '105.43'
    'E6' = '  -6.22%  '
    'D7' = '0.642'
    'E7' = '  +2.42%  '
    'E8' = '  +0.02%  '
    'D9' = '0.649'
    'E9' = '  +7.82%  '
    'D10' = '41.42'
    'E10' = '  -5.31%  '
    'D11' = '0.0938'
    'E11' = '  +2.12%  '
    'D12' = '8.54'
    'E12' = '  -1.19%  '
    'E13' = '  -3.47%  '
    'D14' = '17.02'
    'E14' = '  +13.62%  '
    'E15' = '  +2.17%  '
    'D16' = '2.736.54'
    'E16' = '  +6.68%  '
    'D17' = '2.377.47'
    'E17' = '  +6.61%  '
    'D18' = '43.152.51'
    'E18' = '  +1.53%  '
    'D19' = '7.75'
    'E19' = '  +7.74%  '
    'E20' = '  +2.41%  '
    'D21' = '76.19'
    'E21' = '  +3.55%  '
    'D22' = '275.32'
    'E22' = '  +16.49%  '
    'D23' = '3.39'
    'E23' = '  +0.51%  '
    'D24' = '2.40'
    'E24' = '  +0.68%  '
    'D25' = '9.62'
    'E25' = '  +8.33%  '
    'B26' = 'Dai'
    'C26' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D26' = '0.999'
    'E26' = '  -0.08%  '
    'B27' = 'Cosmos'
    'C27' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D27' = '11.65'
    'E27' = '  +1.73%  '
    'D28' = '22.81'
    'E28' = '  +6.50%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D29' = '2.20'
    'E29' = '  -0.31%  '
    'B30' = 'Monero'
    'C30' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D30' = '175.87'
    'E30' = '  +0.20%  '
    'B31' = 'InjectiveProtocol'
    'C31' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D31' = '37.72'
    'E31' = '  +0.96%  '
    'B32' = 'WEMIXToken'
    'C32' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D32' = '3.21'
    'E32' = '  +2.38%  '
    'D33' = '0.0922'
    'E33' = '  +4.61%  '
    'D34' = '5.86'
    'E34' = '  +3.42%  '
    'D35' = '0.132'
    'E35' = '  +4.77%  '
    'D36' = '4.82'
    'E36' = '  -3.99%  '
    'D37' = '4.10'
    'E37' = '  -2.25%  '
    'D38' = '0.0364'
    'E38' = '  -2.92%  '
    'D39' = '0.106'
    'E39' = '  +2.41%  '
    'D40' = '2.81'
    'E40' = '  +17.85%  '
    'D41' = '1.56'
    'E41' = '  +19.80%  '
    'B42' = 'Aave'
    'C42' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D42' = '124.57'
    'E42' = '  +22.08%  '
    'B43' = 'Algorand'
    'C43' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D43' = '0.230'
    'E43' = '  +0.23%  '
    'D44' = '95.67'
    'E44' = '  +70.66%  '
    'D45' = '68.89'
    'E45' = '  -4.29%  '
    'E46' = '  +0.19%  '
    'D47' = '12.41'
    'E47' = '  +0.84%  '
    'D48' = '9.49'
    'E48' = '  +12.22%  '
    'D49' = '5.59'
    'E49' = '  +3.73%  '
    'D50' = '1.30'
    'E50' = '  +0.70%  '
    'D51' = '0.102'
    'E51' = '  +5.09%  '
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
